$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" table (B15:J24) was re-sorted ascending by period
# (it was previously sorted descending). Update the period labels (E16:E24)
# and the "Valor Mora" amounts (F16:F24) to reflect the new ascending order.

$periods = @("1810", "1811", "1812", "1901", "1902", "1903", "1904", "1905", "1906")
$valores = @(40000, 40000, 40000, 40000, 40000, 40000, 40000, 40000, 22666)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
